$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 2 and row 4 for columns D, L, M, N, O, P, S
# Row 2 becomes what Row 4 was (and vice versa)

# --- Row 2 (new values, taken from old Row 4) ---
$ws.Range("D2").Value = 44923
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7625
$ws.Range("S2").Value = 7625

# --- Row 4 (new values, taken from old Row 2) ---
$ws.Range("D4").Value = 44881
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 11250
$ws.Range("O4").Value = 11250
$ws.Range("P4").Value = 11250
$ws.Range("S4").Value = 11250
